$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.821.19"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "2.344.87"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.11"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.664"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.67"
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.79"
$ws.Range("E11").Value = "  +6.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.03"
$ws.Range("E12").Value = "  +2.87%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.07"
$ws.Range("E15").Value = "  -3.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.900"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "2.347.70"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "43.818.21"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.61"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.48"
$ws.Range("E21").Value = "  -3.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "251.65"
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("B24").Value = "WEMIXToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.79"
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.85"
$ws.Range("E25").Value = "  -3.90%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.37"
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "175.64"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.14"
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0739"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("E34").Value = "  -4.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.74"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.40"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.53"
$ws.Range("E39").Value = "  +18.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0271"
$ws.Range("E40").Value = "  -3.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.34"
$ws.Range("E41").Value = "  +14.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.74"
$ws.Range("E42").Value = "  +2.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.04"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("E44").Value = "  -4.41%  "
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.41"
$ws.Range("E48").Value = "  -3.01%  "
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.38"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("E51").Value = "  +1.75%  "
